$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New date rows (text values, not actual dates) and their Hours/Minutes entries
$rows = @(
    @{ Row = 74; Date = "23/7/2013"; Hours = 1; Minutes = 0 },
    @{ Row = 75; Date = "24/7/2013"; Hours = 1; Minutes = 0 },
    @{ Row = 76; Date = "25/7/2013"; Hours = 1; Minutes = 0 },
    @{ Row = 77; Date = "26/7/2013"; Hours = 0; Minutes = 0 },
    @{ Row = 78; Date = "27/7/2013"; Hours = 0; Minutes = 0 },
    @{ Row = 79; Date = "28/7/2013"; Hours = 0; Minutes = 0 },
    @{ Row = 80; Date = "29/7/2013"; Hours = 2; Minutes = 30 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Hours
    $ws.Cells.Item($r.Row, 3).Value = $r.Minutes
}

# Update the selection to reflect the new active cell (C75) as in the diff
$ws.Range("C75").Select()
